# Updates the crypto price ("D") and 1h volume-change ("E") columns
# on the active sheet to match the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '33.931.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.783.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.30%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.26%  '

$ws.Range("E6").Value = '  -1.34%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.04'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '

$ws.Range("E9").Value = '  +2.99%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0680'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0938'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.039.21'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("E13").Value = '  +2.59%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.786.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.897.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("E16").Value = '  -1.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.79%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0770'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.53%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("E22").Value = '  -2.04%  '

$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("E24").Value = '  -2.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.62%  '

$ws.Range("E26").Value = '  +1.33%  '

$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("E28").Value = '  +0.19%  '

$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("E30").Value = '  +2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0514'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.62'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.56'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.90%  '

$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.393.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("E37").Value = '  -1.33%  '

$ws.Range("E38").Value = '  +1.22%  '

$ws.Range("E39").Value = '  +8.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '79.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("E41").Value = '  +0.60%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.918'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.24%  '

$ws.Range("E44").Value = '  -1.50%  '

$ws.Range("E45").Value = '  +10.93%  '

$ws.Range("E46").Value = '  +3.53%  '

$ws.Range("E47").Value = '  +2.73%  '

$ws.Range("E48").Value = '  +0.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '107.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.941.34'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.24%  '

$ws.Range("E51").Value = '  +0.27%  '
